# "Add files via upload" — replace the participant list (surnames/names) with a
# new roster, add header/footer banner cells with dedicated styles, tidy up
# row heights / selection / sort state / page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe the existing data (18 rows x 2 cols) so stale shared strings don't
#    linger in the rebuilt sharedStrings table.
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
}

# Drop the old sort-state (sortState/sortCondition) left over from the
# previous "Data > Sort" operation.
$ws.Sort.SortFields.Clear()

# ---------------------------------------------------------------------------
# 2. Remove the two now-unused trailing rows (18 & 19) so the sheet is back
#    down to 17 rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()

# ---------------------------------------------------------------------------
# 3. Write the new roster. Row 1 and row 17 are styled "banner" cells with no
#    text; row 2 holds the column headers; rows 3-16 hold the data, which is
#    written in row-major, column-A-before-B order so the shared-strings
#    table is rebuilt in the same sequence as the target file.
# ---------------------------------------------------------------------------
$rows = @(
    @("surnames", "names"),
    @("BILIC", "Marijo"),
    @("DAVIES", "Ben"),
    @("HECTOR", "David"),
    @("HERCOG", "Predrag"),
    @("LEPRI", "Petra"),
    @("MARLIERE", "Fabrice"),
    @("MIHAJLOVIĆ", "Domagoj"),
    @("MRVOS", "Luka"),
    @("NATHALIE", "Bocquet"),
    @("PANAYOTOV", "Nikolay"),
    @("PARVANOVA", "Milena"),
    @("SILHAVY", "Jan"),
    @("VOKOUN", "Miroslav"),
    @("ZOBAJ", "Eugen")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# ---------------------------------------------------------------------------
# 4. Banner-cell formatting.
#    A1: small grey header strip (font size 6, solid #CCCCCC fill, centered
#    vertically, wrapped) with no text.
# ---------------------------------------------------------------------------
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Size = 6
$a1.Interior.Color = 13421772
$a1.WrapText = $true
$a1.VerticalAlignment = -4108
$a1.ClearContents()

#    A17: footer note placeholder (font size 10, centered vertically), no
#    fill, no text.
$a17 = $ws.Cells.Item(17, 1)
$a17.Font.Size = 10
$a17.VerticalAlignment = -4108
$a17.ClearContents()

# ---------------------------------------------------------------------------
# 5. Row heights: rows 5 & 7 get an explicit (custom) 15pt height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15

# ---------------------------------------------------------------------------
# 6. Column width: column A sized to fit ("18" in saved character-width
#    units).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.17

# ---------------------------------------------------------------------------
# 7. Page setup (printer paper size / orientation).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 8. Selection cursor parks on B25 once the edit is done.
# ---------------------------------------------------------------------------
$ws.Range("B25").Select()
